$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.578.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.030.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -7.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.600"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.15"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.370"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.17"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.101"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.325.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.29"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.18"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.762"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.020.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.687.51"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.54"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +10.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0797"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.86"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.70"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.68"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.87"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.38"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0602"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.27"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.97%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.470.09"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0928"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.71"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.71%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.17"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +38.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0203"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.63"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.86"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.93%  "
